$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the A10 cell value from "Assert" to "Then" (bdd syntax fix)
$ws.Range("A10").Value = "Then"

# 2. Rebuild the conditional formatting rules so the differential formats
#    (dxfs) get regenerated the same way the desktop app does whenever a
#    rule is edited through the Conditional Formatting Rules Manager.
$rng = $ws.Range("A1:XFD1048576")
$fcs = $rng.FormatConditions

# remove the existing rules (iterate backwards since collection re-indexes)
for ($i = $fcs.Count; $i -ge 1; $i--) {
    $fcs.Item($i).Delete()
}

# Recreate rule 1: cellIs = "Then" (was "Assert")
$fc1 = $fcs.Add(1, 3, '="Then"')
$fc1.Font.Color = 16777215
$fc1.Interior.Color = 11892015

# Recreate rule 2: cellIs = "When"
$fc2 = $fcs.Add(1, 3, '="When"')
$fc2.Interior.Color = 15652797

# Recreate rule 3: cellIs = "Given a"
$fc3 = $fcs.Add(1, 3, '="Given a"')
$fc3.Interior.Color = 11854022

# Recreate rule 4: cellIs = "Specification"
$fc4 = $fcs.Add(1, 3, '="Specification"')
$fc4.Interior.Color = 11854022

# Recreate rule 5: beginsWith "With Properties"
$fc5 = $fcs.Add(9, [Type]::Missing, [Type]::Missing, [Type]::Missing, "With Properties", 2)
$fc5.Interior.Color = 11389944

# Recreate rule 6: endsWith " table of"
$fc6 = $fcs.Add(9, [Type]::Missing, [Type]::Missing, [Type]::Missing, " table of", 3)
$fc6.Interior.Color = 11389944

# Recreate rule 7: endsWith " of"
$fc7 = $fcs.Add(9, [Type]::Missing, [Type]::Missing, [Type]::Missing, " of", 3)
$fc7.Interior.Color = 10086143

# Recreate rule 8: expression
$fc8 = $fcs.Add(2, [Type]::Missing, 'AND((RIGHT(A1048576, 3) = " of"), A2 = "With Properties")')
$fc8.Font.Color = 16777215
$fc8.Interior.Color = 11892015

# Recreate rule 9: expression
$fc9 = $fcs.Add(2, [Type]::Missing, 'AND(RIGHT(XFD1, 3) = " of", A2 = "With Properties")')
$fc9.Font.Color = 16777215
$fc9.Interior.Color = 11892015

# Recreate rule 10: notContainsBlanks
$fc10 = $fcs.Add(2, [Type]::Missing, 'LEN(TRIM(A1))>0')
$fc10.Font.Color = 16777215
$fc10.Interior.Color = 11892015

# 3. Restore the selected cell shown when the sheet is opened
$ws.Range("A11").Select()
